$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values must be force-formatted as Text before assignment,
# otherwise Excel auto-converts numeric-looking strings (e.g. "213.18") to
# floating point numbers, losing the original text representation.
# NumberFormat is restored to General afterwards to match the source file.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.364.68"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.614.18"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.18"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.250"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.49"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +2.33%  "

$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.68"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.617.06"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.354.75"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.01"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.22"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.62%  "

$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.34"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.87"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.94"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.99%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -2.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.22"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.58"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("E30").Value = "  +4.41%  "

$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("E33").Value = "  -2.24%  "

$ws.Range("E34").Value = "  +2.89%  "

$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.164.48"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +4.91%  "

$ws.Range("E37").Value = "  +3.17%  "

$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.792"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.502"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.786"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.28%  "

$ws.Range("E43").Value = "  +2.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.754.14"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.75"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.53"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.34"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.39%  "

$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.407"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0958"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -10.06%  "

$ws.Range("E51").Value = "  -0.02%  "
